# Update "Run status" sheet: mark additional cases as run ("x") and clear
# cells that no longer apply (cases that are no longer tracked).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run status")

# Cells whose status flips from "r" (running/ref) to "x" (done)
$toMarkX = @(
    "B2","C2","D2","E2","F2",
    "B4","C4","D4","E4","F4",
    "B6","C6","D6","E6","F6",
    "B7","C7","D7","E7","F7",
    "D8","E8","F8",
    "B9","C9","E9","F9",
    "B10","C10","E10",
    "B11","C11","D11","E11","F11",
    "B12","E12","F12",
    "B13","C13","D13",
    "B14","C14","D14","E14","F14",
    "B15","C15","D15","E15",
    "B16","C16","D16","E16","F16",
    "C17","D17","E17","F17",
    "C18","D18","E18"
)

foreach ($addr in $toMarkX) {
    $ws.Range($addr).Value = "x"
}

# Cells that are cleared entirely (case/column combo no longer applicable)
$toClear = @(
    "C5","D5","E5","F5",
    "C8",
    "D9",
    "D10","F10",
    "C12",
    "D12",
    "E13","F13",
    "F15",
    "F18"
)

foreach ($addr in $toClear) {
    $ws.Range($addr).ClearContents()
}

# Update the remembered selection on the sheet (was F19, now B19)
$ws.Range("B19").Select()
